$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.749.60"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "3.751.07"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'621.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "'180.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").Value = "3.750.21"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("E11").Value = "  -5.35%  "
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").Value = "'40.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").Value = "4.368.67"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").Value = "3.750.30"
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").Value = "69.826.78"
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'7.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").Value = "'16.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").Value = "'506.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.99%  "
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("D26").Value = "'13.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.25%  "
$ws.Range("D27").Value = "'11.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("E28").Value = "  +25.73%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("D33").Value = "'31.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").Value = "'0.115"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +3.98%  "
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("E40").Value = "  -5.10%  "
$ws.Range("D41").Value = "'50.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.64%  "
$ws.Range("D42").Value = "'45.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").Value = "'426.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.90%  "
$ws.Range("B45").Value = "Cosmos"
$ws.Range("C45").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D45").Value = "'8.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").Value = "3.006.17"
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").Value = "'27.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.10%  "
$ws.Range("D50").Value = "'137.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("E51").Value = "  +1.26%  "
